$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New headers on row 9 - written in the order that reproduces the
#    target shared-string table ordering (Base, Tercera, Primera, Segunga).
# ---------------------------------------------------------------------
$ws.Range("G9").Value = "Base"
$ws.Range("P9").Value = "Tercera Iteración"
$ws.Range("J9").Value = "Primera Iteración "
$ws.Range("M9").Value = "Segunga Iteración"

# ---------------------------------------------------------------------
# 2. Four side-by-side mini tables (rows 10-15).
# ---------------------------------------------------------------------

# --- Base (G:H) ---
$ws.Cells.Item(10,7).Value = 1
$ws.Cells.Item(10,8).Value = 10
$ws.Cells.Item(11,7).Value = 1
$ws.Cells.Item(11,8).Value = 10
$ws.Cells.Item(12,7).Value = 2
$ws.Cells.Item(12,8).Value = 110
$ws.Cells.Item(13,7).Value = 2
$ws.Cells.Item(13,8).Value = 110
$ws.Cells.Item(14,7).Value = 4
$ws.Cells.Item(14,8).Value = 20
$ws.Cells.Item(15,7).Value = 4
$ws.Cells.Item(15,8).Value = 20

# --- Primera Iteración (J:K) ---
$ws.Cells.Item(10,10).Value = 4
$ws.Cells.Item(10,11).Value = 20
$ws.Cells.Item(11,10).Value = 2
$ws.Cells.Item(11,11).Value = 110
$ws.Cells.Item(12,10).Value = 2
$ws.Cells.Item(12,11).Value = 110
$ws.Cells.Item(13,10).Value = 2
$ws.Cells.Item(13,11).Value = 110
$ws.Cells.Item(14,10).Value = 4
$ws.Cells.Item(14,11).Value = 20
$ws.Cells.Item(15,10).Value = 4
$ws.Cells.Item(15,11).Value = 20

# --- Segunga Iteración (M:N) ---
$ws.Cells.Item(10,13).Value = 2
$ws.Cells.Item(10,14).Value = 110
$ws.Cells.Item(11,13).Value = 1
$ws.Cells.Item(11,14).Value = 10
$ws.Cells.Item(12,13).Value = 2
$ws.Cells.Item(12,14).Value = 110
$ws.Cells.Item(13,13).Value = 2
$ws.Cells.Item(13,14).Value = 110
$ws.Cells.Item(14,13).Value = 4
$ws.Cells.Item(14,14).Value = 20
$ws.Cells.Item(15,13).Value = 4
$ws.Cells.Item(15,14).Value = 20

# --- Tercera Iteración (P:Q) ---
$ws.Cells.Item(10,16).Value = 3
$ws.Cells.Item(10,17).Value = 30
# P11:Q11 intentionally left blank (styled only, no values)
$ws.Cells.Item(12,16).Value = 2
$ws.Cells.Item(12,17).Value = 110
$ws.Cells.Item(13,16).Value = 2
$ws.Cells.Item(13,17).Value = 110
$ws.Cells.Item(14,16).Value = 4
$ws.Cells.Item(14,17).Value = 20
$ws.Cells.Item(15,16).Value = 4
$ws.Cells.Item(15,17).Value = 20

# ---------------------------------------------------------------------
# 3. Formatting. Themed fills are applied first (they share a color
#    table entry regardless of order), then the plain RGB fill and the
#    grid border, which keeps the generated style table minimal.
# ---------------------------------------------------------------------
$ws.Range("C6:H6").Interior.ThemeColor = 8
$ws.Range("C6:H6").Interior.TintAndShade = 0.59999389629810485

$ws.Range("J10:K15").Interior.ThemeColor = 9
$ws.Range("J10:K15").Interior.TintAndShade = 0.59999389629810485

$ws.Range("M10:N15").Interior.ThemeColor = 8
$ws.Range("M10:N15").Interior.TintAndShade = 0.59999389629810485

$ws.Range("P10:Q15").Interior.ThemeColor = 10
$ws.Range("P10:Q15").Interior.TintAndShade = 0.59999389629810485

$ws.Range("G10:H15").Interior.Color = 65535

$gridRange = $ws.Range("C5:I7")
$gridRange.Borders.LineStyle = 1
$gridRange.Borders.Weight = 2

# ---------------------------------------------------------------------
# 4. Totals row - same bold style as the existing H16 total.
# ---------------------------------------------------------------------
$ws.Range("K16").Formula = "=SUM(K10:K15)"
$ws.Range("K16").Font.Bold = $true

$ws.Range("N16").Formula = "=SUM(N10:N15)"
$ws.Range("N16").Font.Bold = $true

$ws.Range("Q16").Formula = "=SUM(Q10:Q15)"
$ws.Range("Q16").Font.Bold = $true

# ---------------------------------------------------------------------
# 5. Column I best-fit width and final selection.
# ---------------------------------------------------------------------
$ws.Columns("I").ColumnWidth = 13.28515625

$ws.Range("K8").Select() | Out-Null
